$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D (rows 2-19 and 21-22) to "N", leaving D20 unchanged ("Y")
foreach ($r in 2..19) {
    $ws.Cells.Item($r, 4).Value = "N"
}
foreach ($r in 21..22) {
    $ws.Cells.Item($r, 4).Value = "N"
}

# Update the selection to D20
$ws.Range("D20").Select()
